$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" ("72.420.08")
Set-TextCell $ws "E2" ("  +1.04%  ")

# Row 3
Set-TextCell $ws "D3" ("4.043.55")
Set-TextCell $ws "E3" ("  +0.61%  ")

# Row 4
Set-TextCell $ws "E4" ("  -0.16%  ")

# Row 5
Set-TextCell $ws "D5" ("542.30")
Set-TextCell $ws "E5" ("  +1.48%  ")

# Row 6
Set-TextCell $ws "D6" ("152.01")
Set-TextCell $ws "E6" ("  -0.43%  ")

# Row 7
Set-TextCell $ws "D7" ("4.035.31")
Set-TextCell $ws "E7" ("  +0.60%  ")

# Row 8
Set-TextCell $ws "D8" ("0.697")
Set-TextCell $ws "E8" ("  +0.84%  ")

# Row 9
Set-TextCell $ws "D9" ("0.999")
Set-TextCell $ws "E9" ("  -0.08%  ")

# Row 10
Set-TextCell $ws "E10" ("  +0.21%  ")

# Row 11
Set-TextCell $ws "D11" ("0.173")
Set-TextCell $ws "E11" ("  +0.47%  ")

# Row 12
Set-TextCell $ws "D12" ("53.72")
Set-TextCell $ws "E12" ("  +11.78%  ")

# Row 13
Set-TextCell $ws "D13" ("0.0000332")
Set-TextCell $ws "E13" ("  +1.74%  ")

# Row 14
Set-TextCell $ws "D14" ("10.95")
Set-TextCell $ws "E14" ("  +1.56%  ")

# Row 15
Set-TextCell $ws "D15" ("4.685.61")
Set-TextCell $ws "E15" ("  +0.49%  ")

# Row 16
Set-TextCell $ws "D16" ("4.040.75")
Set-TextCell $ws "E16" ("  +0.23%  ")

# Row 17
Set-TextCell $ws "D17" ("14.37")
Set-TextCell $ws "E17" ("  +1.58%  ")

# Row 18
Set-TextCell $ws "D18" ("20.75")
Set-TextCell $ws "E18" ("  +1.13%  ")

# Row 19
Set-TextCell $ws "D19" ("1.21")
Set-TextCell $ws "E19" ("  +1.09%  ")

# Row 20
Set-TextCell $ws "E20" ("  -0.69%  ")

# Row 21
Set-TextCell $ws "D21" ("72.354.26")
Set-TextCell $ws "E21" ("  +0.96%  ")

# Row 22
Set-TextCell $ws "D22" ("450.05")
Set-TextCell $ws "E22" ("  +4.30%  ")

# Row 23
Set-TextCell $ws "D23" ("98.02")
Set-TextCell $ws "E23" ("  -0.83%  ")

# Row 24
Set-TextCell $ws "D24" ("3.54")
Set-TextCell $ws "E24" ("  +0.28%  ")

# Row 25
Set-TextCell $ws "D25" ("4.28")
Set-TextCell $ws "E25" ("  +2.03%  ")

# Row 26
Set-TextCell $ws "D26" ("14.64")
Set-TextCell $ws "E26" ("  +0.85%  ")

# Row 27
Set-TextCell $ws "D27" ("4.25")
Set-TextCell $ws "E27" ("  +13.91%  ")

# Row 28
Set-TextCell $ws "D28" ("11.29")
Set-TextCell $ws "E28" ("  +1.58%  ")

# Row 29
Set-TextCell $ws "D29" ("10.82")
Set-TextCell $ws "E29" ("  +0.09%  ")

# Row 30
Set-TextCell $ws "D30" ("5.96")
Set-TextCell $ws "E30" ("  +2.19%  ")

# Row 31
Set-TextCell $ws "D31" ("37.18")
Set-TextCell $ws "E31" ("  +0.84%  ")

# Row 32
Set-TextCell $ws "D32" ("7.97")
Set-TextCell $ws "E32" ("  +14.24%  ")

# Row 33
Set-TextCell $ws "E33" ("  +3.04%  ")

# Row 34
Set-TextCell $ws "D34" ("13.63")
Set-TextCell $ws "E34" ("  +1.42%  ")

# Row 35
Set-TextCell $ws "D35" ("49.12")
Set-TextCell $ws "E35" ("  +16.31%  ")

# Row 36
Set-TextCell $ws "D36" ("681.36")
Set-TextCell $ws "E36" ("  -2.34%  ")

# Row 37
Set-TextCell $ws "E37" ("  +0.07%  ")

# Row 38
Set-TextCell $ws "E38" ("  +5.66%  ")

# Row 39
Set-TextCell $ws "D39" ("0.0" + [string][char]0x2083 + "0894")
Set-TextCell $ws "E39" ("  +8.33%  ")

# Row 40
Set-TextCell $ws "D40" ("3.43")
Set-TextCell $ws "E40" ("  -1.83%  ")

# Row 41
Set-TextCell $ws "D41" ("0.149")
Set-TextCell $ws "E41" ("  -4.13%  ")

# Row 42
Set-TextCell $ws "D42" ("3.39")
Set-TextCell $ws "E42" ("  -2.16%  ")

# Row 43
Set-TextCell $ws "D43" ("11.19")
Set-TextCell $ws "E43" ("  +17.41%  ")

# Row 44
Set-TextCell $ws "D44" ("0.999")
Set-TextCell $ws "E44" ("  +0.02%  ")

# Row 45
Set-TextCell $ws "B45" ("VeChain")
Set-TextCell $ws "C45" ("https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet")
Set-TextCell $ws "D45" ("0.0496")
Set-TextCell $ws "E45" ("  +1.81%  ")

# Row 46
Set-TextCell $ws "B46" ("FirstDigitalUSD")
Set-TextCell $ws "C46" ("https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd")
Set-TextCell $ws "D46" ("0.999")
Set-TextCell $ws "E46" ("  -0.05%  ")

# Row 47
Set-TextCell $ws "E47" ("  +0.70%  ")

# Row 48
Set-TextCell $ws "D48" ("2.71")
Set-TextCell $ws "E48" ("  +2.65%  ")

# Row 49
Set-TextCell $ws "D49" ("3.12")
Set-TextCell $ws "E49" ("  +3.44%  ")

# Row 50
Set-TextCell $ws "B50" ("LidoDAOToken")
Set-TextCell $ws "C50" ("https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo")
Set-TextCell $ws "D50" ("3.52")
Set-TextCell $ws "E50" ("  +5.89%  ")

# Row 51
Set-TextCell $ws "B51" ("ApeXProtocol")
Set-TextCell $ws "C51" ("https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex")
Set-TextCell $ws "D51" ("3.30")
Set-TextCell $ws "E51" ("  -2.74%  ")
